# SO Test Plan, Migration from QARSF to AutomationOrg
#
# The CashReceipt sheet gains a new "CustomerID" column (with sample
# value "a5B41000000PRNXEA4") inserted right before the existing
# "Customer" column, and the existing "Cust-Dollar WF1 (8)" label is
# renamed to "Cust-Dollar WF1".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CashReceipt")

# Insert a new blank column before column D ("Customer"); this shifts the
# old D/E/F columns (Customer / CustomerPaymentAmount / CustomerPaymentMethod)
# one place to the right, to E/F/G.
$ws.Columns.Item(4).Insert()

# Header for the new column, formatted like the other bold header cells
# (e.g. "Deposit Amount Bank" in B1 / "Customer" now in E1).
$ws.Cells.Item(1,4).Value = "CustomerID"
$ws.Cells.Item(1,2).Copy() | Out-Null
$ws.Cells.Item(1,4).PasteSpecial(-4122) | Out-Null

# Sample value for the new column, rendered in Arial (rather than the
# workbook default Calibri).
$ws.Cells.Item(2,4).Value = "a5B41000000PRNXEA4"
$tmpStyle = $wb.Styles.Add("TmpCustomerIdStyle")
$tmpStyle.Font.Name = "Arial"
$ws.Cells.Item(2,4).Style = "TmpCustomerIdStyle"
$tmpStyle.Delete()

# The "Customer" data cell (now E2) gets the trailing " (8)" dropped from
# its label.
$ws.Cells.Item(2,5).Value = "Cust-Dollar WF1"

# Matches the author's last-saved selection on the sheet.
$ws.Range("I9").Select() | Out-Null
